$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Unified_table")
$ws.Range("A43").Value = 23
